# Completes the PEXTRA (grib table 126) comment strings for several
# "identified missing CMIP6 requested variables" rows, per commit:
# "The comment of the PEXTRA table 126 variables in the (pre) identified
#  missing file has been completed #450."
#
# Every edit below only touches column H ("comment") text for the rows
# whose comment needed the extra "grib 126.xx ... part of MFPxxx" /
# "To be implemented: ..." annotation. Rows that merely duplicate the
# same variable further down the sheet (and therefore already carry an
# identical comment string elsewhere) are updated too, since in the
# source workbook they shared the same underlying text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H8").Value  = 'To be implemented:  grib 126.30  part of MFPPHY   For Greenland this is the same as above sftgif. We do not have Antarctic ice sheet.'
$ws.Range("H13").Value = 'grib 126.105                                                                   part of MFP3D        Available in IFS: T-tendency from convection : grib 128.105'
$ws.Range("H15").Value = 'COSP grib 126.46   CVEXTR2(7)=''ISCCP_MEANALBCLD''       part of MFPPHY'
$ws.Range("H16").Value = 'COSP grib 126.42   CVEXTR2(3)=''CALIPSO_HCC''                part of MFPPHY'
$ws.Range("H17").Value = 'COSP grib 126.40   CVEXTR2(1)=''CALIPSO_LCC''                part of MFPPHY'
$ws.Range("H18").Value = 'COSP grib 126.41   CVEXTR2(2)=''CALIPSO_MCC''               part of MFPPHY'
$ws.Range("H19").Value = 'COSP grib 126.43   CVEXTR2(4)=''CALIPSO_TCC''                part of MFPPHY'
$ws.Range("H20").Value = 'COSP grib 126.44   CVEXTR2(5)=''ISCCP_TOTALCLD''       part of MFPPHY'
$ws.Range("H21").Value = 'COSP grib 126.45   CVEXTR2(6)=''ISCCP_MEANPTOP''           part of MFPPHY'
$ws.Range("H22").Value = 'Grib 126.94 + 126.99 + 126.106 + 126.110       part of MFP3D        Adding all the q-tendencies, thus: grib 128.94 + 128.99 + 128.106 + 128.110.  Alternatively, in IFS: just estimating the delta q per month. So far no direct grib code for the totoal q-tendency found'
$ws.Range("H23").Value = 'grib 126.106                                                                   part of MFP3D        Available in IFS: q-tendency from convection: grib 128.106'
$ws.Range("H24").Value = 'grib 126.99 + 126.106 + 126.110                           part of MFP3D        Adding all the q-tendencies without advection, thus: grib 128.99 + 128.106 + 128.110.'
$ws.Range("H25").Value = 'grib 126.105                                                                   part of MFP3D        Available in IFS: T-tendency from convection : grib 128.105'
$ws.Range("H26").Value = 'grib 126.95                                                                      part of MFP3D        Available in IFS: T-tendency from radiation: grib 128.95'
$ws.Range("H28").Value = 'To be implemented:  grib 126.34  part of MFPPHY   Available in PISM. This is the ice sheet mask (in fraction) defined in the ice sheet model grid'
$ws.Range("H32").Value = 'To be implemented:  grib 126.34  part of MFPPHY   Available in PISM. This is the ice sheet mask (in fraction) defined in the ice sheet model grid'
$ws.Range("H34").Value = 'Grib 126.20 / 126.22        part of MFP3D        In namelist.ifs.cloudact+diag.sh  CVEXTRA(1)=''CDNC'' which is a PEXTRA variable.'
$ws.Range("H35").Value = 'grib 126.73                          part of MFPPHY    Available from double radiation call in IFS. PEXTRA issue #403   aerosol free'
$ws.Range("H36").Value = 'grib 126.72                          part of MFPPHY    Available from double radiation call in IFS. PEXTRA issue #403   aerosol free'
$ws.Range("H37").Value = 'grib 128.212-126.069     part of MFPPHY    Available from double radiation call in IFS. PEXTRA issue #403   aerosol free'
$ws.Range("H47").Value = 'COSP grib 126.46   CVEXTR2(7)=''ISCCP_MEANALBCLD''       part of MFPPHY'
$ws.Range("H48").Value = 'COSP grib 126.42   CVEXTR2(3)=''CALIPSO_HCC''                part of MFPPHY'
$ws.Range("H49").Value = 'COSP grib 126.40   CVEXTR2(1)=''CALIPSO_LCC''                part of MFPPHY'
$ws.Range("H50").Value = 'COSP grib 126.41   CVEXTR2(2)=''CALIPSO_MCC''               part of MFPPHY'
$ws.Range("H51").Value = 'COSP grib 126.43   CVEXTR2(4)=''CALIPSO_TCC''                part of MFPPHY'
$ws.Range("H52").Value = 'COSP grib 126.44   CVEXTR2(5)=''ISCCP_TOTALCLD''       part of MFPPHY'
$ws.Range("H53").Value = 'COSP grib 126.45   CVEXTR2(6)=''ISCCP_MEANPTOP''           part of MFPPHY'
$ws.Range("H57").Value = 'COSP grib 126.42   CVEXTR2(3)=''CALIPSO_HCC''                part of MFPPHY'
$ws.Range("H58").Value = 'COSP grib 126.40   CVEXTR2(1)=''CALIPSO_LCC''                part of MFPPHY'
$ws.Range("H59").Value = 'COSP grib 126.41   CVEXTR2(2)=''CALIPSO_MCC''               part of MFPPHY'
$ws.Range("H60").Value = 'COSP grib 126.43   CVEXTR2(4)=''CALIPSO_TCC''                part of MFPPHY'
$ws.Range("H63").Value = 'Grib 126.94 + 126.99 + 126.106 + 126.110       part of MFP3D        Adding all the q-tendencies, thus: grib 128.94 + 128.99 + 128.106 + 128.110.  Alternatively, in IFS: just estimating the delta q per month. So far no direct grib code for the totoal q-tendency found'
$ws.Range("H65").Value = 'To be implemented:  grib 126.32  part of MFPPHY   This is the land ice mask and will be an extra variable in IFS (thomas: via PEXTRA?)'
$ws.Range("H66").Value = 'To be implemented:  grib 126.30  part of MFPPHY   For Greenland this is the same as above sftgif. We do not have Antarctic ice sheet.'
$ws.Range("H67").Value = 'To be implemented:  grib 126.31  part of MFPPHY   Not available in IFS. Although it could be calculated from tile fractions and written out as extra output'
$ws.Range("H71").Value = 'Grib 126.94 + 126.99 + 126.106 + 126.110       part of MFP3D        Adding all the q-tendencies, thus: grib 128.94 + 128.99 + 128.106 + 128.110.  Alternatively, in IFS: just estimating the delta q per month. So far no direct grib code for the totoal q-tendency found'
$ws.Range("H72").Value = 'grib 126.106                                                                   part of MFP3D        Available in IFS: q-tendency from convection: grib 128.106'
$ws.Range("H73").Value = 'grib 126.99 + 126.106 + 126.110                           part of MFP3D        Adding all the q-tendencies without advection, thus: grib 128.99 + 128.106 + 128.110.'
$ws.Range("H74").Value = 'grib 126.105                                                                   part of MFP3D        Available in IFS: T-tendency from convection : grib 128.105'
$ws.Range("H75").Value = 'grib 126.95                                                                      part of MFP3D        Available in IFS: T-tendency from radiation: grib 128.95'

# Note: row 65 (sftgif, table LImon) gets the completed comment above.
# Row 81 is the duplicate sftgif entry further down (table IyrGre) and in
# the source workbook keeps the original, not-yet-annotated comment text
# ("This is the land ice mask and will be an extra variable in IFS
# (thomas: via PEXTRA?)") -- so it is intentionally left untouched here.

# Best-effort: the source diff also nudges the saved scroll position of
# the sheet view (topLeftCell A66 -> A55); reproduce that intent so the
# view starts a bit higher up, if the host supports it.
try {
  $excel.ActiveWindow.ScrollRow = 55
} catch {
}
